# Add a new "Hungary" worksheet to the workbook, based on the existing
# "Slovakia" sheet (same layout/styles), with its own market name and
# Jira reference, and update the view/selection state of both sheets to
# reflect Hungary becoming the newly active tab.

$wb = $excel.ActiveWorkbook

$slovakia = $wb.Worksheets.Item("Slovakia")

# Duplicate the Slovakia sheet (placed immediately after it) and rename it.
$slovakia.Copy($null, $slovakia) | Out-Null
$hungary = $wb.Worksheets.Item($wb.Worksheets.Count)
$hungary.Name = "Hungary"

# Update the market name / reference cells for the new Hungary sheet.
$hungary.Range("B2").Value = "Hungary Market"
$hungary.Range("B4").Value = "NGC-4308/T3594/T3619"

# Move the Slovakia sheet's selection off of its old cell (it is no longer
# the active tab).
$slovakia.Activate() | Out-Null
$slovakia.Cells.Select() | Out-Null

# Hungary becomes the active sheet/tab, with its own cell selection.
$hungary.Activate() | Out-Null
$hungary.Range("A13").Select() | Out-Null
